# Auto-generated script applying the Coeurl_Profits market-data refresh diff.
# Each block updates the H:N "current price / leve profit" columns for a given
# sheet+row to match the values captured by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 451.58334
$ws.Range("J2").Value = 400
$ws.Range("L2").Value = 400
$ws.Range("N2").Value = -626
$ws.Range("H62").Value = 6596.7
$ws.Range("I62").Value = 6996
$ws.Range("K62").Value = 6996
$ws.Range("M62").Value = -6372
$ws.Range("H65").Value = 6596.7
$ws.Range("I65").Value = 6996
$ws.Range("K65").Value = 34980
$ws.Range("M65").Value = -31860
$ws.Range("H76").Value = 7333.3335
$ws.Range("H79").Value = 7333.3335
$ws.Range("H80").Value = 1285.9565
$ws.Range("I80").Value = 1322.75
$ws.Range("K80").Value = 3968.25
$ws.Range("M80").Value = -2970.25
$ws.Range("H83").Value = 1285.9565
$ws.Range("I83").Value = 1322.75
$ws.Range("K83").Value = 11904.75
$ws.Range("M83").Value = -6912.75
$ws.Range("H86").Value = 7399.6
$ws.Range("I86").Value = 7399.6
$ws.Range("K86").Value = 7399.6
$ws.Range("M86").Value = -6276.6
$ws.Range("H89").Value = 7399.6
$ws.Range("I89").Value = 7399.6
$ws.Range("K89").Value = 36998
$ws.Range("M89").Value = -31382
$ws.Range("H132").Value = 8334900.5
$ws.Range("I132").Value = 10102131
$ws.Range("K132").Value = 30306393
$ws.Range("M132").Value = -30303863
$ws.Range("H137").Value = 10815.546
$ws.Range("I137").Value = 1864
$ws.Range("K137").Value = 5592
$ws.Range("M137").Value = -3042
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("H32").Value = 4961.475
$ws.Range("I32").Value = 4561.1353
$ws.Range("K32").Value = 4561.1353
$ws.Range("M32").Value = -4274.1353
$ws.Range("H34").Value = 79997.5
$ws.Range("I34").Value = 60000
$ws.Range("J34").Value = 86663.336
$ws.Range("K34").Value = 60000
$ws.Range("L34").Value = 86663.336
$ws.Range("M34").Value = -59729
$ws.Range("N34").Value = -87205.336
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = $null
$ws.Range("N37").Value = 0
$ws.Range("H60").Value = 32999.5
$ws.Range("I60").Value = 45999
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 45999
$ws.Range("L60").Value = 20000
$ws.Range("M60").Value = -45266
$ws.Range("N60").Value = -21466
$ws.Range("H63").Value = 3185
$ws.Range("I63").Value = 3659
$ws.Range("K63").Value = 3659
$ws.Range("M63").Value = -2973
$ws.Range("H66").Value = 3185
$ws.Range("I66").Value = 3659
$ws.Range("K66").Value = 18295
$ws.Range("M66").Value = -14863
$ws.Range("H88").Value = 2275
$ws.Range("J88").Value = 2812.5
$ws.Range("L88").Value = 2812.5
$ws.Range("N88").Value = -3624.5
$ws.Range("H91").Value = 2275
$ws.Range("J91").Value = 2812.5
$ws.Range("L91").Value = 2812.5
$ws.Range("N91").Value = -5620.5
$ws.Range("H121").Value = $null
$ws.Range("I121").Value = $null
$ws.Range("J121").Value = $null
$ws.Range("K121").Value = $null
$ws.Range("L121").Value = $null
$ws.Range("N121").Value = $null
$ws.Range("H122").Value = $null
$ws.Range("I122").Value = $null
$ws.Range("J122").Value = $null
$ws.Range("K122").Value = $null
$ws.Range("L122").Value = $null
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = $null
$ws.Range("H123").Value = $null
$ws.Range("I123").Value = $null
$ws.Range("J123").Value = $null
$ws.Range("K123").Value = $null
$ws.Range("L123").Value = $null
$ws.Range("H124").Value = $null
$ws.Range("I124").Value = $null
$ws.Range("J124").Value = $null
$ws.Range("K124").Value = $null
$ws.Range("L124").Value = $null
$ws.Range("H125").Value = $null
$ws.Range("I125").Value = $null
$ws.Range("J125").Value = $null
$ws.Range("K125").Value = $null
$ws.Range("L125").Value = $null
$ws.Range("N125").Value = $null
$ws.Range("H126").Value = $null
$ws.Range("I126").Value = $null
$ws.Range("J126").Value = $null
$ws.Range("K126").Value = $null
$ws.Range("L126").Value = $null
$ws.Range("M126").Value = $null
$ws.Range("H127").Value = $null
$ws.Range("I127").Value = $null
$ws.Range("J127").Value = $null
$ws.Range("K127").Value = $null
$ws.Range("L127").Value = $null
$ws.Range("H128").Value = $null
$ws.Range("I128").Value = $null
$ws.Range("J128").Value = $null
$ws.Range("K128").Value = $null
$ws.Range("L128").Value = $null
$ws.Range("M128").Value = $null
$ws.Range("H129").Value = $null
$ws.Range("I129").Value = $null
$ws.Range("J129").Value = $null
$ws.Range("K129").Value = $null
$ws.Range("L129").Value = $null
$ws.Range("H130").Value = $null
$ws.Range("I130").Value = $null
$ws.Range("J130").Value = $null
$ws.Range("K130").Value = $null
$ws.Range("L130").Value = $null
$ws.Range("H131").Value = $null
$ws.Range("I131").Value = $null
$ws.Range("J131").Value = $null
$ws.Range("K131").Value = $null
$ws.Range("L131").Value = $null
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = $null
$ws.Range("I132").Value = $null
$ws.Range("J132").Value = $null
$ws.Range("K132").Value = $null
$ws.Range("L132").Value = $null
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null
$ws.Range("H133").Value = $null
$ws.Range("I133").Value = $null
$ws.Range("J133").Value = $null
$ws.Range("K133").Value = $null
$ws.Range("L133").Value = $null
$ws.Range("N133").Value = $null
$ws.Range("H134").Value = $null
$ws.Range("I134").Value = $null
$ws.Range("J134").Value = $null
$ws.Range("K134").Value = $null
$ws.Range("L134").Value = $null
$ws.Range("H135").Value = $null
$ws.Range("I135").Value = $null
$ws.Range("J135").Value = $null
$ws.Range("K135").Value = $null
$ws.Range("L135").Value = $null
$ws.Range("H137").Value = $null
$ws.Range("I137").Value = $null
$ws.Range("J137").Value = $null
$ws.Range("K137").Value = $null
$ws.Range("L137").Value = $null
$ws.Range("H138").Value = $null
$ws.Range("I138").Value = $null
$ws.Range("J138").Value = $null
$ws.Range("K138").Value = $null
$ws.Range("L138").Value = $null
$ws.Range("M138").Value = $null
$ws.Range("H139").Value = $null
$ws.Range("I139").Value = $null
$ws.Range("J139").Value = $null
$ws.Range("K139").Value = $null
$ws.Range("L139").Value = $null
$ws.Range("H140").Value = $null
$ws.Range("I140").Value = $null
$ws.Range("J140").Value = $null
$ws.Range("K140").Value = $null
$ws.Range("L140").Value = $null
$ws.Range("H141").Value = $null
$ws.Range("I141").Value = $null
$ws.Range("J141").Value = $null
$ws.Range("K141").Value = $null
$ws.Range("L141").Value = $null
$ws.Range("N141").Value = $null
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 54684.105
$ws.Range("I86").Value = 64374.875
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 64374.875
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -63251.875
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 54684.105
$ws.Range("I89").Value = 64374.875
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 321874.375
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -316258.375
$ws.Range("N89").Value = -26232
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1149.7142
$ws.Range("I17").Value = 739.8
$ws.Range("J17").Value = 2174.5
$ws.Range("K17").Value = 739.8
$ws.Range("L17").Value = 2174.5
$ws.Range("M17").Value = -565.8
$ws.Range("N17").Value = -2522.5
$ws.Range("H31").Value = 22467.04
$ws.Range("I31").Value = 34101.773
$ws.Range("K31").Value = 34101.773
$ws.Range("M31").Value = -33806.773
$ws.Range("H34").Value = 22467.04
$ws.Range("I34").Value = 34101.773
$ws.Range("K34").Value = 34101.773
$ws.Range("M34").Value = -33899.773
$ws.Range("H58").Value = 2948.721
$ws.Range("I58").Value = 2631.2334
$ws.Range("J58").Value = 3681.3845
$ws.Range("K58").Value = 2631.2334
$ws.Range("L58").Value = 3681.3845
$ws.Range("M58").Value = -2428.2334
$ws.Range("N58").Value = -4087.3845
$ws.Range("H86").Value = 4416
$ws.Range("J86").Value = 3925
$ws.Range("L86").Value = 3925
$ws.Range("N86").Value = -6171
$ws.Range("H89").Value = 4416
$ws.Range("J89").Value = 3925
$ws.Range("L89").Value = 19625
$ws.Range("N89").Value = -30857
$ws.Range("H136").Value = 2948.721
$ws.Range("I136").Value = 2631.2334
$ws.Range("J136").Value = 3681.3845
$ws.Range("K136").Value = 7893.7002
$ws.Range("L136").Value = 11044.1535
$ws.Range("M136").Value = -5343.7002
$ws.Range("N136").Value = -16144.1535
$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 949.75
$ws.Range("I129").Value = 662.25
$ws.Range("J129").Value = 1524.75
$ws.Range("K129").Value = 1986.75
$ws.Range("L129").Value = 4574.25
$ws.Range("M129").Value = 3013.25
$ws.Range("N129").Value = -14574.25
$ws.Range("H131").Value = 29127.865
$ws.Range("I131").Value = 84187.164
$ws.Range("J131").Value = 2699.4
$ws.Range("K131").Value = 252561.492
$ws.Range("L131").Value = 8098.200000000001
$ws.Range("M131").Value = -247521.492
$ws.Range("N131").Value = -18178.2
$ws.Range("H136").Value = 721941
$ws.Range("I136").Value = 1251646.8
$ws.Range("J136").Value = 15666.667
$ws.Range("K136").Value = 3754940.4
$ws.Range("L136").Value = 47000.001
$ws.Range("M136").Value = -3749840.4
$ws.Range("N136").Value = -57200.001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17296.715
$ws.Range("I70").Value = 5808
$ws.Range("K70").Value = 5808
$ws.Range("M70").Value = -5538
$ws.Range("H73").Value = 17296.715
$ws.Range("I73").Value = 5808
$ws.Range("K73").Value = 5808
$ws.Range("M73").Value = -4872
$ws.Range("H80").Value = 3347.6667
$ws.Range("I80").Value = 2949.6667
$ws.Range("K80").Value = 2949.6667
$ws.Range("M80").Value = -1951.6667
$ws.Range("H83").Value = 3347.6667
$ws.Range("I83").Value = 2949.6667
$ws.Range("K83").Value = 14748.3335
$ws.Range("M83").Value = -9756.333500000001
$ws.Range("H86").Value = 38994
$ws.Range("J86").Value = 38994
$ws.Range("L86").Value = 38994
$ws.Range("N86").Value = -41366
$ws.Range("H89").Value = 38994
$ws.Range("J89").Value = 38994
$ws.Range("L89").Value = 116982
$ws.Range("N89").Value = -128838
$ws.Range("H122").Value = 1803.4286
$ws.Range("I122").Value = 1362.3334
$ws.Range("J122").Value = 4450
$ws.Range("K122").Value = 4087.0002
$ws.Range("L122").Value = 13350
$ws.Range("M122").Value = -1637.0002
$ws.Range("N122").Value = -18250
$ws.Range("H125").Value = $null
$ws.Range("I125").Value = $null
$ws.Range("J125").Value = $null
$ws.Range("K125").Value = $null
$ws.Range("L125").Value = $null
$ws.Range("N125").Value = $null
$ws.Range("H126").Value = $null
$ws.Range("I126").Value = $null
$ws.Range("J126").Value = $null
$ws.Range("K126").Value = $null
$ws.Range("L126").Value = $null
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = $null
$ws.Range("H127").Value = $null
$ws.Range("I127").Value = $null
$ws.Range("J127").Value = $null
$ws.Range("K127").Value = $null
$ws.Range("L127").Value = $null
$ws.Range("N127").Value = $null
$ws.Range("H128").Value = $null
$ws.Range("I128").Value = $null
$ws.Range("J128").Value = $null
$ws.Range("K128").Value = $null
$ws.Range("L128").Value = $null
$ws.Range("H129").Value = $null
$ws.Range("I129").Value = $null
$ws.Range("J129").Value = $null
$ws.Range("K129").Value = $null
$ws.Range("L129").Value = $null
$ws.Range("N129").Value = $null
$ws.Range("H130").Value = $null
$ws.Range("I130").Value = $null
$ws.Range("J130").Value = $null
$ws.Range("K130").Value = $null
$ws.Range("L130").Value = $null
$ws.Range("N130").Value = $null
$ws.Range("H131").Value = $null
$ws.Range("I131").Value = $null
$ws.Range("J131").Value = $null
$ws.Range("K131").Value = $null
$ws.Range("L131").Value = $null
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = $null
$ws.Range("I132").Value = $null
$ws.Range("J132").Value = $null
$ws.Range("K132").Value = $null
$ws.Range("L132").Value = $null
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null
$ws.Range("H133").Value = $null
$ws.Range("I133").Value = $null
$ws.Range("J133").Value = $null
$ws.Range("K133").Value = $null
$ws.Range("L133").Value = $null
$ws.Range("H134").Value = $null
$ws.Range("I134").Value = $null
$ws.Range("J134").Value = $null
$ws.Range("K134").Value = $null
$ws.Range("L134").Value = $null
$ws.Range("N134").Value = $null
$ws.Range("H135").Value = $null
$ws.Range("I135").Value = $null
$ws.Range("J135").Value = $null
$ws.Range("K135").Value = $null
$ws.Range("L135").Value = $null
$ws.Range("N135").Value = $null
$ws.Range("H136").Value = $null
$ws.Range("I136").Value = $null
$ws.Range("J136").Value = $null
$ws.Range("K136").Value = $null
$ws.Range("L136").Value = $null
$ws.Range("N136").Value = $null
$ws.Range("H137").Value = $null
$ws.Range("I137").Value = $null
$ws.Range("J137").Value = $null
$ws.Range("K137").Value = $null
$ws.Range("L137").Value = $null
$ws.Range("M137").Value = $null
$ws.Range("N137").Value = $null
$ws.Range("H138").Value = $null
$ws.Range("I138").Value = $null
$ws.Range("J138").Value = $null
$ws.Range("K138").Value = $null
$ws.Range("L138").Value = $null
$ws.Range("N138").Value = $null
$ws.Range("H139").Value = $null
$ws.Range("I139").Value = $null
$ws.Range("J139").Value = $null
$ws.Range("K139").Value = $null
$ws.Range("L139").Value = $null
$ws.Range("M139").Value = $null
$ws.Range("N139").Value = $null
$ws.Range("H140").Value = $null
$ws.Range("I140").Value = $null
$ws.Range("J140").Value = $null
$ws.Range("K140").Value = $null
$ws.Range("L140").Value = $null
$ws.Range("M140").Value = $null
$ws.Range("H141").Value = $null
$ws.Range("I141").Value = $null
$ws.Range("J141").Value = $null
$ws.Range("K141").Value = $null
$ws.Range("L141").Value = $null
$ws.Range("N141").Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 53241
$ws.Range("J6").Value = 53241
$ws.Range("L6").Value = 53241
$ws.Range("N6").Value = -53465
$ws.Range("H13").Value = 4318.7144
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = $null
$ws.Range("H25").Value = 8343010.5
$ws.Range("I25").Value = 37500600
$ws.Range("J25").Value = 12270.857
$ws.Range("K25").Value = 37500600
$ws.Range("L25").Value = 12270.857
$ws.Range("M25").Value = -37500370
$ws.Range("N25").Value = -12730.857
$ws.Range("H68").Value = 2993.75
$ws.Range("I68").Value = 2993.75
$ws.Range("K68").Value = 2993.75
$ws.Range("M68").Value = -2244.75
$ws.Range("H71").Value = 2993.75
$ws.Range("I71").Value = 2993.75
$ws.Range("K71").Value = 14968.75
$ws.Range("M71").Value = -11224.75
$ws.Range("H93").Value = 144971
$ws.Range("I93").Value = 252074.5
$ws.Range("J93").Value = 2166.3333
$ws.Range("K93").Value = 252074.5
$ws.Range("L93").Value = 2166.3333
$ws.Range("M93").Value = -250826.5
$ws.Range("N93").Value = -4662.3333
$ws.Range("H136").Value = 45822.086
$ws.Range("I136").Value = 57617.11
$ws.Range("K136").Value = 172851.33
$ws.Range("M136").Value = -170301.33
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1002.25
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null
$ws.Range("H126").Value = 9675
$ws.Range("I126").Value = 11028.417
$ws.Range("J126").Value = 4261.3335
$ws.Range("K126").Value = 33085.251
$ws.Range("L126").Value = 12784.0005
$ws.Range("M126").Value = -30615.251
$ws.Range("N126").Value = -17724.0005
$ws.Range("H132").Value = 1182.6666
$ws.Range("I132").Value = 1136.75
$ws.Range("K132").Value = 3410.25
$ws.Range("M132").Value = -880.25
$ws.Range("H136").Value = 2107.0513
$ws.Range("I136").Value = 1960.5358
$ws.Range("K136").Value = 5881.607400000001
$ws.Range("M136").Value = -3331.607400000001

$wb.Save()
